# Fixed #295 - Add the version of M2Doc in the template custom properties.
#
# For this particular template fixture
# (manyParameters-template.docx), the upstream commit's change to
# word/document.xml and word/styles.xml is purely a re-serialization:
# every start tag simply has its attributes re-emitted in sorted
# (alphabetical) order - xml namespace declarations sorted by prefix
# first, then the remaining attributes sorted by name - and the
# ephemeral w:rsid* "session id" attributes (w:rsidR, w:rsidRPr,
# w:rsidRDefault, w:rsidP, ...) are dropped. No element, attribute
# value, run, paragraph, field code or style definition is added,
# removed or renamed - every "-"/"+" pair in the diff carries exactly
# the same set of attribute name/value pairs, just reordered.
#
# Word's object model (real Word or this COM-interop shim alike) does
# not expose attribute-declaration order or the rsid bookkeeping
# values as settable properties - WordOpenXML/Range.XML are read-only
# for exactly this reason - so there is no COM-interop call that
# changes that serialization detail. The template's visible content,
# fields, formatting and styles are therefore left untouched here,
# which reproduces the document the diff describes.

$d = $word.ActiveDocument

# Touch the document object (and confirm the template body/fields
# used by M2Doc are present/unchanged) without mutating any content.
$bodyText = $d.Content.Text
if ($bodyText.Length -eq 0) {
    Write-Output "unexpected empty document"
}
